$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (header) is unchanged: Scientific name | Author | Rank | Parent name

# Row 2 (Eukaryota kingdom) is unchanged: A2=Eukaryota, C2=Kingdom

# Row 3: was "Unidentified coccoid eukaryots" / Phylum / Eukaryota -> becomes the Prokaryota kingdom row
$ws.Range("A3").Value = "Prokaryota"
$ws.Range("C3").Value = "Kingdom"
$ws.Range("D3").ClearContents()

# Row 4 previously held "Unidentified flagellated eukaryots" - it is removed entirely (blank separator row)
$ws.Rows("4:4").ClearContents()

# Block: Flagellates incertae sedis chain (rows 5-10), child of Eukaryota
$ws.Range("A5").Value = "Flagellates phylum incertae sedis"
$ws.Range("C5").Value = "Phylum"
$ws.Range("D5").Value = "Eukaryota"

$ws.Range("A6").Value = "Flagellates classes incertae sedis"
$ws.Range("C6").Value = "Class"
$ws.Range("D6").Value = "Flagellates phylum incertae sedis"

$ws.Range("A7").Value = "Flagellates ordines incertae sedis"
$ws.Range("C7").Value = "Order"
$ws.Range("D7").Value = "Flagellates classes incertae sedis"

$ws.Range("A8").Value = "Flagellates families incertae sedis"
$ws.Range("C8").Value = "Family"
$ws.Range("D8").Value = "Flagellates ordines incertae sedis"

$ws.Range("A9").Value = "Flagellates genera incertae sedis"
$ws.Range("C9").Value = "Genus"
$ws.Range("D9").Value = "Flagellates families incertae sedis"

$ws.Range("A10").Value = "Flagellates species incertae sedis"
$ws.Range("C10").Value = "Species"
$ws.Range("D10").Value = "Flagellates genera incertae sedis"

# Row 11 left blank (separator row)

# Block: Unicells incertae sedis chain (rows 12-18), starting with its own kingdom row
$ws.Range("A12").Value = "Unicells kingdom incertae sedis"
$ws.Range("C12").Value = "Kingdom"

$ws.Range("A13").Value = "Unicells phylum incertae sedis"
$ws.Range("C13").Value = "Phylum"
$ws.Range("D13").Value = "Unicells kingdom incertae sedis"

$ws.Range("A14").Value = "Unicells classes incertae sedis"
$ws.Range("C14").Value = "Class"
$ws.Range("D14").Value = "Unicells phylum incertae sedis"

$ws.Range("A15").Value = "Unicells ordines incertae sedis"
$ws.Range("C15").Value = "Order"
$ws.Range("D15").Value = "Unicells classes incertae sedis"

$ws.Range("A16").Value = "Unicells families incertae sedis"
$ws.Range("C16").Value = "Family"
$ws.Range("D16").Value = "Unicells ordines incertae sedis"

$ws.Range("A17").Value = "Unicells genera incertae sedis"
$ws.Range("C17").Value = "Genus"
$ws.Range("D17").Value = "Unicells families incertae sedis"

$ws.Range("A18").Value = "Unicells species incertae sedis"
$ws.Range("C18").Value = "Species"
$ws.Range("D18").Value = "Unicells genera incertae sedis"

# Row 19 left blank (separator row)

# Block: Eukarotic picoplankton incertae sedis chain (rows 20-25), child of Eukaryota
$ws.Range("A20").Value = "Eukarotic picoplankton phylum incertae sedis"
$ws.Range("C20").Value = "Phylum"
$ws.Range("D20").Value = "Eukaryota"

$ws.Range("A21").Value = "Eukarotic picoplankton classes incertae sedis"
$ws.Range("C21").Value = "Class"
$ws.Range("D21").Value = "Eukarotic picoplankton phylum incertae sedis"

$ws.Range("A22").Value = "Eukarotic picoplankton ordines incertae sedis"
$ws.Range("C22").Value = "Order"
$ws.Range("D22").Value = "Eukarotic picoplankton classes incertae sedis"

$ws.Range("A23").Value = "Eukarotic picoplankton families incertae sedis"
$ws.Range("C23").Value = "Family"
$ws.Range("D23").Value = "Eukarotic picoplankton ordines incertae sedis"

$ws.Range("A24").Value = "Eukarotic picoplankton genera incertae sedis"
$ws.Range("C24").Value = "Genus"
$ws.Range("D24").Value = "Eukarotic picoplankton families incertae sedis"

$ws.Range("A25").Value = "Eukarotic picoplankton species incertae sedis"
$ws.Range("C25").Value = "Species"
$ws.Range("D25").Value = "Eukarotic picoplankton genera incertae sedis"
